# Refresh the "Data" sheet of the USDSOFRCSA_USD curve workbook.
#
# The old single 1Y SOFROIS/OIS row is dropped and replaced by ten new
# FUTURE quotes (tenors 3M..15M plus 0M/1M/3M), while the remaining OIS
# tenors (2Y..30Y) are kept as-is but shift down by 9 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Drop the old row 3 ("1Y" / SOFROIS / OIS / 0.038348) - it is not present
# in the refreshed table; rows 4:9 (2Y..30Y) move up to 3:8.
$ws.Range("A3").EntireRow.Delete()

# Make room for the 10 replacement FUTURE rows; this pushes the former
# 2Y..30Y block (now at rows 3:8) back down to rows 13:18.
$ws.Range("A3:A12").EntireRow.Insert()

$futures = @(
    @("3M",  "SQZ25", "FUTURE", 96.20999999999999),
    @("5M",  "SQF26", "FUTURE", 96.30500000000001),
    @("5M",  "SQG26", "FUTURE", 96.395),
    @("6M",  "SQH26", "FUTURE", 96.43000000000001),
    @("9M",  "SQM26", "FUTURE", 96.65000000000001),
    @("12M", "SQU26", "FUTURE", 96.8),
    @("15M", "SQZ26", "FUTURE", 96.875),
    @("0M",  "SQU25", "FUTURE", 95.9025),
    @("1M",  "SQV25", "FUTURE", 96.01000000000001),
    @("3M",  "SQX25", "FUTURE", 96.125)
)

$r = 3
foreach ($row in $futures) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
